$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text (Question -> Questions, Answer -> Answers)
$ws.Range("A1").Value = "Questions"
$ws.Range("B1").Value = "Answers"

# Set active cell/selection to B6 as in the diff
$ws.Range("B6").Select()
